# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values for the e9d423f6-... row (row 4) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-23 06:46:43"
$wsZhCn.Range("H4").Value = "2016-03-23 06:47:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-23 06:46:51"
$wsDeDe.Range("H4").Value = "2016-03-23 06:47:41"
